# Update cryptos list values: refreshed prices/volume percentages, and
# re-ranked rows 17/18 (WrappedEther now above Polkadot) and 34/35/36
# (EnergySwap, NEARProtocol, Aptos) per the latest coinranking.com pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.716.73'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '3.095.07'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''516.72'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '''143.60'
$ws.Range("E6").Value = '  +3.51%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").Value = '''0.437'
$ws.Range("E8").Value = '  +1.18%  '
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("E10").Value = '  -0.02%  '
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").Value = '3.609.12'
$ws.Range("E12").Value = '  +2.48%  '
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("D14").Value = '''25.86'
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("D15").Value = '''0.0000166'
$ws.Range("E15").Value = '  -1.85%  '
$ws.Range("D16").Value = '57.819.39'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.098.20'
$ws.Range("E17").Value = '  +2.76%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '''6.16'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").Value = '''13.13'
$ws.Range("E19").Value = '  -1.32%  '
$ws.Range("D20").Value = '''8.23'
$ws.Range("E20").Value = '  +1.65%  '
$ws.Range("D21").Value = '''337.99'
$ws.Range("E21").Value = '  +3.34%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '''0.503'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '''65.77'
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("D25").Value = '''0.173'
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").Value = '0.0₃0938'
$ws.Range("E27").Value = '  +5.93%  '
$ws.Range("D28").Value = '''6.50'
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("D29").Value = '''7.12'
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").Value = '''20.94'
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("E32").Value = '  -2.43%  '
$ws.Range("D33").Value = '''154.82'
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("B34").Value = 'EnergySwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D34").Value = '''28.43'
$ws.Range("E34").Value = '  +12.42%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '''4.55'
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '''5.93'
$ws.Range("E36").Value = '  +1.96%  '
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").Value = '3.133.18'
$ws.Range("E39").Value = '  +2.59%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").Value = '''0.674'
$ws.Range("E42").Value = '  +2.29%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '2.285.49'
$ws.Range("E44").Value = '  +5.20%  '
$ws.Range("D45").Value = '''0.0255'
$ws.Range("E45").Value = '  +4.87%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '''20.44'
$ws.Range("E47").Value = '  +4.02%  '
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("D49").Value = '''5.90'
$ws.Range("E49").Value = '  -4.22%  '
$ws.Range("D50").Value = '''0.0879'
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("E51").Value = '  +3.09%  '
